$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - same inline-string / bold-bordered-centered style as
# the existing header cells (copy the format from H1, which carries style index 1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells (rows 2 and 3), no special style - matches C2:H3 etc.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 6
